$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "district"
$ws.Range("E2").Value = "Aligarh"
